# ---------------------------------------------------------------------------
# Applies the "Updated App Technical Document" revision to
# Android_App_Data_Collection - V2.docx
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaRange($text) {
    # Returns a fresh Range(start,end) spanning exactly the text content
    # (no trailing paragraph mark) of the paragraph whose text matches $text.
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -eq ($text + "`r")) {
            return $d.Range($p.Range.Start, $p.Range.End - 1)
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. "Launch count" bullet: append a new run after "user opened the app "
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("user opened the app ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $tgt = $d.Range($r.Start, $r.End)
    $tgt.Collapse(0)
    $tgt.InsertAfter("from the beginning of the year")
}

# ---------------------------------------------------------------------------
# 2. Bold note: "The data collection period is one year." -> add clause
# ---------------------------------------------------------------------------
$rng = Get-ParaRange("Note: The data collection period is one year.")
if ($rng -ne $null) {
    $xml = "<w:p $W><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`">Note: </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>The dat</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>a collection period is one year, starting from the beginning of the year.</w:t></w:r></w:p>"
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3. Add _GoBack bookmark at the end of the BootReceiver.java bullet
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("turns their phone back on.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $tgt = $d.Range($r.End - 1, $r.End)
    $d.Bookmarks.Add("_GoBack", $tgt)
}

# ---------------------------------------------------------------------------
# 4. MyService.java bullet: remove comma before "and buffering"
# ---------------------------------------------------------------------------
$rng = Get-ParaRange("This file contains the logic for data collection, and buffering.")
if ($rng -ne $null) {
    $xml = "<w:p $W><w:r><w:t>This file contains</w:t></w:r><w:r><w:t xml:space=`"preserve`"> the logic for data collection </w:t></w:r><w:r><w:t>and buffering</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 5. Insert three new bullets after the MyService.java description
#    (numId=4, same as that bullet's list)
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$anchorPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "This file contains the logic for data collection and buffering.`r") {
        $anchorPara = $p
        break
    }
}
if ($anchorPara -ne $null) {
    $anchorPara.Range.InsertParagraphAfter()

    $paras = $d.Paragraphs
    $newPara = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -eq "`r" -and $p.Range.Start -eq ($anchorPara.Range.End)) {
            $newPara = $p
            break
        }
    }
    if ($newPara -ne $null) {
        $xml1 = "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"4`"/></w:numPr></w:pPr>" +
                "<w:r><w:t xml:space=`"preserve`">Both Start and </w:t></w:r>" +
                "<w:r><w:t>End tim</w:t></w:r>" +
                "<w:r><w:t>e</w:t></w:r>" +
                "<w:r><w:t xml:space=`"preserve`"> are calculated usi</w:t></w:r>" +
                "<w:r><w:t xml:space=`"preserve`">ng the </w:t></w:r>" +
                "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/>" +
                "<w:r><w:t>getLastTimeUsed</w:t></w:r>" +
                "<w:proofErr w:type=`"spellEnd`"/>" +
                "<w:r><w:t>(</w:t></w:r>" +
                "<w:proofErr w:type=`"gramEnd`"/>" +
                "<w:r><w:t>) method</w:t></w:r>" +
                "<w:r><w:t xml:space=`"preserve`"> provided by the API</w:t></w:r>" +
                "<w:r><w:t>.</w:t></w:r></w:p>"
        $rng0 = $d.Range($newPara.Range.Start, $newPara.Range.Start)
        $rng0.InsertXML($xml1)
    }

    # find the paragraph we just filled in, then insert the next one after it
    $paras = $d.Paragraphs
    $p1 = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -eq "Both Start and End time are calculated using the getLastTimeUsed() method provided by the API.`r") {
            $p1 = $p
            break
        }
    }
    if ($p1 -ne $null) {
        $p1.Range.InsertParagraphAfter()

        $paras = $d.Paragraphs
        $newPara2 = $null
        for ($i = 1; $i -le $paras.Count; $i++) {
            $p = $paras.Item($i)
            if ($p.Range.Text -eq "`r" -and $p.Range.Start -eq ($p1.Range.End)) {
                $newPara2 = $p
                break
            }
        }
        if ($newPara2 -ne $null) {
            $xml2 = "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"4`"/></w:numPr></w:pPr>" +
                    "<w:r><w:t>Usage statistics</w:t></w:r>" +
                    "<w:r><w:t xml:space=`"preserve`"> of an app are checked within a 5 second interval to monitor any change and update the data accordingly. This is done in </w:t></w:r>" +
                    "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/>" +
                    "<w:r><w:t>CollectData</w:t></w:r>" +
                    "<w:proofErr w:type=`"spellEnd`"/>" +
                    "<w:r><w:t>(</w:t></w:r>" +
                    "<w:proofErr w:type=`"gramEnd`"/>" +
                    "<w:r><w:t>).</w:t></w:r></w:p>"
            $rng1 = $d.Range($newPara2.Range.Start, $newPara2.Range.Start)
            $rng1.InsertXML($xml2)
        }

        $paras = $d.Paragraphs
        $p2 = $null
        for ($i = 1; $i -le $paras.Count; $i++) {
            $p = $paras.Item($i)
            if ($p.Range.Text -eq "Usage statistics of an app are checked within a 5 second interval to monitor any change and update the data accordingly. This is done in CollectData().`r") {
                $p2 = $p
                break
            }
        }
        if ($p2 -ne $null) {
            $p2.Range.InsertParagraphAfter()

            $paras = $d.Paragraphs
            $newPara3 = $null
            for ($i = 1; $i -le $paras.Count; $i++) {
                $p = $paras.Item($i)
                if ($p.Range.Text -eq "`r" -and $p.Range.Start -eq ($p2.Range.End)) {
                    $newPara3 = $p
                    break
                }
            }
            if ($newPara3 -ne $null) {
                $xml3 = "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"4`"/></w:numPr></w:pPr>" +
                        "<w:r><w:t xml:space=`"preserve`">An array called </w:t></w:r>" +
                        "<w:proofErr w:type=`"spellStart`"/>" +
                        "<w:r><w:t>statsList</w:t></w:r>" +
                        "<w:proofErr w:type=`"spellEnd`"/>" +
                        "<w:r><w:t xml:space=`"preserve`"> contains the usage statistic data collected from apps</w:t></w:r>" +
                        "<w:r><w:t>.</w:t></w:r></w:p>"
                $rng2 = $d.Range($newPara3.Range.Start, $newPara3.Range.Start)
                $rng2.InsertXML($xml3)
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 6. Ustats.java bullet: "functions" -> "methods" (split run) + drop old
#    _GoBack bookmark, merging "UsageStatsM" + "anager" into one run
# ---------------------------------------------------------------------------
$rng = Get-ParaRange("This file has functions that implement the UsageStatsManager API and allows the app to collect data.")
if ($rng -ne $null) {
    $xml = "<w:p $W><w:r><w:t>This file has method</w:t></w:r>" +
           "<w:r><w:t xml:space=`"preserve`">s that implement the </w:t></w:r>" +
           "<w:proofErr w:type=`"spellStart`"/>" +
           "<w:r><w:t>UsageStatsManager</w:t></w:r>" +
           "<w:proofErr w:type=`"spellEnd`"/>" +
           "<w:r><w:t xml:space=`"preserve`"> API and allows the app to collect data</w:t></w:r>" +
           "<w:r><w:t>.</w:t></w:r></w:p>"
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 7. Compare.java heading: insert a leading lastRenderedPageBreak + br run
# ---------------------------------------------------------------------------
$rng = Get-ParaRange("Compare.java")
if ($rng -ne $null) {
    $xml = "<w:p $W><w:r><w:lastRenderedPageBreak/><w:br/></w:r><w:r><w:t>Compare.java</w:t></w:r></w:p>"
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 8. Compare.java bullet: split run (no visible text change) + new bullet
#    "A list of these objects is sent to the server." (numId=5)
# ---------------------------------------------------------------------------
$rng = Get-ParaRange("This file represents the app usage statistics we" + [char]8217 + "re collecting, and is used in MyService.java.")
if ($rng -ne $null) {
    $xml = "<w:p $W><w:r><w:t>This file represents the app usage statistics we" + [char]8217 + "re collecting</w:t></w:r><w:r><w:t>, and is used in MyService.java.</w:t></w:r></w:p>"
    $rng.InsertXML($xml)
}

$paras = $d.Paragraphs
$anchorPara2 = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq ("This file represents the app usage statistics we" + [char]8217 + "re collecting, and is used in MyService.java.`r")) {
        $anchorPara2 = $p
        break
    }
}
if ($anchorPara2 -ne $null) {
    $anchorPara2.Range.InsertParagraphAfter()

    $paras = $d.Paragraphs
    $newPara4 = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -eq "`r" -and $p.Range.Start -eq ($anchorPara2.Range.End)) {
            $newPara4 = $p
            break
        }
    }
    if ($newPara4 -ne $null) {
        $newPara4.Range.InsertAfter("A list of these objects is sent to the server.")
    }
}

Write-Output "edit complete"
